$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 92947
$ws.Range("E2").Value = -2702
$ws.Range("F2").Value = -2702
$ws.Range("G2").Value = -3991
$ws.Range("H2").Value = -4405
$ws.Range("I2").Value = -4540
$ws.Range("J2").Value = 135
$ws.Range("K2").Value = 105967
$ws.Range("L2").Value = 61315
$ws.Range("M2").Value = 44653
$ws.Range("N2").Value = 42281
$ws.Range("O2").Value = 2371
$ws.Range("P2").Value = 2185
$ws.Range("Q2").Value = 317
$ws.Range("R2").Value = -3486
$ws.Range("S2").Value = 1249
$ws.Range("T2").Value = 1101
$ws.Range("U2").Value = -784
$ws.Range("V2").Value = 19050
$ws.Range("W2").Value = -2.91
$ws.Range("X2").Value = -4.74
$ws.Range("Y2").Value = -10.15
$ws.Range("Z2").Value = -4.12
$ws.Range("AA2").Value = 137.31
$ws.Range("AB2").Value = 1828.22
$ws.Range("AC2").Value = -11762
$ws.Range("AD2").Value = -5.59
$ws.Range("AE2").Value = 109537
$ws.Range("AF2").Value = 0.6
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 0.15
$ws.Range("AI2").Value = -0.89
$ws.Range("AJ2").Value = 34800000
# Row 3
$ws.Range("D3").Value = 95137
$ws.Range("E3").Value = 2718
$ws.Range("F3").Value = 2718
$ws.Range("G3").Value = 4362
$ws.Range("H3").Value = 2170
$ws.Range("I3").Value = 2068
$ws.Range("J3").Value = 102
$ws.Range("K3").Value = 120649
$ws.Range("L3").Value = 72591
$ws.Range("M3").Value = 48058
$ws.Range("N3").Value = 43448
$ws.Range("O3").Value = 4610
$ws.Range("P3").Value = 2185
$ws.Range("Q3").Value = 3849
$ws.Range("R3").Value = -1020
$ws.Range("S3").Value = 5107
$ws.Range("T3").Value = 1955
$ws.Range("U3").Value = 1894
$ws.Range("V3").Value = 29547
$ws.Range("W3").Value = 2.86
$ws.Range("X3").Value = 2.28
$ws.Range("Y3").Value = 4.82
$ws.Range("Z3").Value = 1.92
$ws.Range("AA3").Value = 151.05
$ws.Range("AB3").Value = 1911.89
$ws.Range("AC3").Value = 5357
$ws.Range("AD3").Value = 12.53
$ws.Range("AE3").Value = 112559
$ws.Range("AF3").Value = 0.6
$ws.Range("AG3").Value = 300
$ws.Range("AH3").Value = 0.45
$ws.Range("AI3").Value = 5.69
$ws.Range("AJ3").Value = 34800000
# Row 4
$ws.Range("D4").Value = 98538
$ws.Range("E4").Value = 4194
$ws.Range("F4").Value = 4194
$ws.Range("G4").Value = 4364
$ws.Range("H4").Value = 2932
$ws.Range("I4").Value = 2653
$ws.Range("J4").Value = 279
$ws.Range("K4").Value = 123915
$ws.Range("L4").Value = 72461
$ws.Range("M4").Value = 51454
$ws.Range("N4").Value = 46109
$ws.Range("O4").Value = 5344
$ws.Range("P4").Value = 2185
$ws.Range("Q4").Value = 1470
$ws.Range("R4").Value = -6994
$ws.Range("S4").Value = -463
$ws.Range("T4").Value = 2541
$ws.Range("U4").Value = -1071
$ws.Range("V4").Value = 28277
$ws.Range("W4").Value = 4.26
$ws.Range("X4").Value = 2.97
$ws.Range("Y4").Value = 5.92
$ws.Range("Z4").Value = 2.4
$ws.Range("AA4").Value = 140.83
$ws.Range("AB4").Value = 2039.97
$ws.Range("AC4").Value = 6873
$ws.Range("AD4").Value = 12.67
$ws.Range("AE4").Value = 119454
$ws.Range("AF4").Value = 0.73
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 0.34
$ws.Range("AI4").Value = 4.44
$ws.Range("AJ4").Value = 34800000
# Row 5
$ws.Range("D5").Value = 123355
$ws.Range("E5").Value = 5459
$ws.Range("F5").Value = 5459
$ws.Range("G5").Value = 6753
$ws.Range("H5").Value = 5080
$ws.Range("I5").Value = 4905
$ws.Range("J5").Value = 175
$ws.Range("K5").Value = 134025
$ws.Range("L5").Value = 77081
$ws.Range("M5").Value = 56943
$ws.Range("N5").Value = 50643
$ws.Range("O5").Value = 6301
$ws.Range("P5").Value = 2185
$ws.Range("Q5").Value = 4376
$ws.Range("R5").Value = -2374
$ws.Range("S5").Value = 1539
$ws.Range("T5").Value = 2360
$ws.Range("U5").Value = 2017
$ws.Range("V5").Value = 31545
$ws.Range("W5").Value = 4.42
$ws.Range("X5").Value = 4.12
$ws.Range("Y5").Value = 10.14
$ws.Range("Z5").Value = 3.94
$ws.Range("AA5").Value = 135.37
$ws.Range("AB5").Value = 2259.89
$ws.Range("AC5").Value = 12707
$ws.Range("AD5").Value = 6.48
$ws.Range("AE5").Value = 131198
$ws.Range("AF5").Value = 0.63
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1.21
$ws.Range("AI5").Value = 7.91
$ws.Range("AJ5").Value = 34800000
# Row 6
$ws.Range("D6").Value = 109845
$ws.Range("E6").Value = 8454
$ws.Range("F6").Value = 8454
$ws.Range("G6").Value = 8944
$ws.Range("H6").Value = 6781
$ws.Range("I6").Value = 6464
$ws.Range("K6").Value = 128337
$ws.Range("L6").Value = 67835
$ws.Range("M6").Value = 60502
$ws.Range("N6").Value = 53900
$ws.Range("P6").Value = 2185
$ws.Range("Q6").Value = 11045
$ws.Range("R6").Value = -3808
$ws.Range("S6").Value = -5256
$ws.Range("T6").Value = 1485
$ws.Range("U6").Value = 9560
$ws.Range("V6").Value = 26594
$ws.Range("W6").Value = 7.7
$ws.Range("X6").Value = 6.17
$ws.Range("Y6").Value = 12.37
$ws.Range("Z6").Value = 5.17
$ws.Range("AA6").Value = 112.12
$ws.Range("AB6").Value = 2440.21
$ws.Range("AC6").Value = 16746
$ws.Range("AD6").Value = 6.12
$ws.Range("AE6").Value = 139638
$ws.Range("AF6").Value = 0.73
$ws.Range("AG6").Value = 1700
$ws.Range("AH6").Value = 1.66
$ws.Range("AI6").Value = 10.18
$ws.Range("AJ6").Value = 34800000
# Row 7
$ws.Range("D7").Value = 95188
$ws.Range("E7").Value = 10161
$ws.Range("G7").Value = 10064
$ws.Range("H7").Value = 7555
$ws.Range("I7").Value = 7037
$ws.Range("K7").Value = 130211
$ws.Range("L7").Value = 63225
$ws.Range("M7").Value = 66986
$ws.Range("N7").Value = 59758
$ws.Range("P7").Value = 2188
$ws.Range("Q7").Value = 6027
$ws.Range("R7").Value = -1908
$ws.Range("S7").Value = -2705
$ws.Range("T7").Value = 1660
$ws.Range("U7").Value = 4128
$ws.Range("W7").Value = 10.67
$ws.Range("X7").Value = 7.94
$ws.Range("Y7").Value = 12.38
$ws.Range("Z7").Value = 5.84
$ws.Range("AA7").Value = 94.39
$ws.Range("AC7").Value = 18230
$ws.Range("AD7").Value = 4.67
$ws.Range("AE7").Value = 154813
$ws.Range("AF7").Value = 0.55
$ws.Range("AG7").Value = 1816
$ws.Range("AH7").Value = 2.13
$ws.Range("AI7").Value = 8.98
# Row 8
$ws.Range("D8").Value = 106620
$ws.Range("E8").Value = 10357
$ws.Range("G8").Value = 10520
$ws.Range("H8").Value = 7824
$ws.Range("I8").Value = 7318
$ws.Range("K8").Value = 138257
$ws.Range("L8").Value = 64544
$ws.Range("M8").Value = 73714
$ws.Range("N8").Value = 65980
$ws.Range("P8").Value = 2188
$ws.Range("Q8").Value = 7679
$ws.Range("R8").Value = -5738
$ws.Range("S8").Value = -1298
$ws.Range("T8").Value = 2871
$ws.Range("U8").Value = 2784
$ws.Range("W8").Value = 9.710000000000001
$ws.Range("X8").Value = 7.34
$ws.Range("Y8").Value = 11.64
$ws.Range("Z8").Value = 5.83
$ws.Range("AA8").Value = 87.56
$ws.Range("AC8").Value = 18958
$ws.Range("AD8").Value = 4.27
$ws.Range("AE8").Value = 170934
$ws.Range("AF8").Value = 0.47
$ws.Range("AG8").Value = 1819
$ws.Range("AH8").Value = 2.25
$ws.Range("AI8").Value = 8.65
# Row 9
$ws.Range("D9").Value = 109068
$ws.Range("E9").Value = 10172
$ws.Range("G9").Value = 10423
$ws.Range("H9").Value = 7692
$ws.Range("I9").Value = 7236
$ws.Range("K9").Value = 145128
$ws.Range("L9").Value = 64490
$ws.Range("M9").Value = 80638
$ws.Range("N9").Value = 72450
$ws.Range("P9").Value = 2188
$ws.Range("Q9").Value = 8274
$ws.Range("R9").Value = -4912
$ws.Range("S9").Value = -1317
$ws.Range("T9").Value = 2218
$ws.Range("U9").Value = 4664
$ws.Range("W9").Value = 9.33
$ws.Range("X9").Value = 7.05
$ws.Range("Y9").Value = 10.46
$ws.Range("Z9").Value = 5.43
$ws.Range("AA9").Value = 79.98
$ws.Range("AC9").Value = 18747
$ws.Range("AD9").Value = 4.32
$ws.Range("AE9").Value = 187693
$ws.Range("AF9").Value = 0.43
$ws.Range("AG9").Value = 1829
$ws.Range("AH9").Value = 2.26
$ws.Range("AI9").Value = 8.789999999999999
